# Book log workbook edit:
#  - Convert the remaining book's ISBN-10 to ISBN-13
#  - Remove the two other book rows (rows 3 & 4), keeping only the
#    "Morrie" record, with updated catalog details.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 3 and 4 entirely (their content, and shift rows below up).
$ws.Rows("3:4").Delete()

# Update the remaining data row (row 2) with the new book details.
$ws.Range("A2").Value = 9780385318792      # ISBN-13 (converted from ISBN-10 385318790)
$ws.Range("B2").Value = "Schwartz, Morrie"
$ws.Range("C2").Value = "Delta"
$ws.Range("D2").Value = "Morrie: In His Own Words"
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 56
